$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Add a new quarter sheet "2022-Q4" right after "总计" (before the
#    current "2022-Q1" sheet). Copying the existing "2022-Q1" sheet
#    carries over its column layout/number formats/borders so the new
#    sheet looks exactly like its siblings.
# ------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Copy($q1)
$newSheet = $wb.Worksheets.Item("2022-Q1 (2)")
$newSheet.Name = "2022-Q4"

# The copied sheet only has 2 data rows (rows 2:3); we need 6 (rows
# 2:7). Extend the formatting of the existing data rows down so the
# extra rows pick up the same per-column formatting.
$newSheet.Range("A2:H3").Copy()
$newSheet.Range("A4:H7").PasteSpecial(-4122)

# Fund-code / numeric-text columns must stay textual (leading zeros in
# fund codes like "004707" must survive), matching the source data's
# use of literal text for every column except A (index) and H (rank).
$newSheet.Range("B2:G7").NumberFormat = "@"

$newSheetRows = @(
  @("004707", "景顺长城睿成灵活配置混合A", "5.12", "36.44", "0.77", "0.0394", 7),
  @("014466", "工银行业优选混合A", "1.09", "62.72", "3.43", "0.0374", 6),
  @("014467", "工银行业优选混合C", "0.54", "62.72", "3.43", "0.0185", 6),
  @("004719", "景顺长城睿成灵活配置混合C", "1.12", "36.44", "0.77", "0.0086", 7),
  @("165524", "信诚中证智能家居指数（LOF）A", "0.37", "91.28", "1.00", "0.0037", 6),
  @("013084", "信诚中证智能家居指数（LOF）C", "0.14", "91.28", "1.00", "0.0014", 6)
)

for ($i = 0; $i -lt $newSheetRows.Count; $i++) {
  $r = 2 + $i
  $row = $newSheetRows[$i]
  $newSheet.Range("A$r").Value = $i
  $newSheet.Range("B$r").Value = $row[0]
  $newSheet.Range("C$r").Value = $row[1]
  $newSheet.Range("D$r").Value = $row[2]
  $newSheet.Range("E$r").Value = $row[3]
  $newSheet.Range("F$r").Value = $row[4]
  $newSheet.Range("G$r").Value = $row[5]
  $newSheet.Range("H$r").Value = $row[6]
}

# ------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row for 2022-Q4
#    right under the header row, pushing every other quarter down by
#    one row.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows(2).Insert()

# Re-apply the per-column formatting (index column keeps its style,
# the rest stay unstyled) by copying it down from the row that used
# to be row 2 and now lives at row 3.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B3:D3").Copy()
$summary.Range("B2:D2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 6
$summary.Range("D2").Value = 0.11

# Renumber the running index in column A for the rows that shifted
# down (they used to read 0..5, now need to read 1..6).
for ($r = 3; $r -le 8; $r++) {
  $summary.Range("A$r").Value = $r - 2
}
